# Added Test Data for Croatia Market
#
# Duplicate the "Spain" sheet (the last sheet) into a new "Croatia" sheet,
# insert the extra "MZX Communicator" row that the other 16-row country
# sheets (e.g. Germany) already have, and set the market-specific values.

$wb = $excel.ActiveWorkbook

$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Copy([Type]::Missing, $spain)

$croatia = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Name = "Croatia"

# Make room for the "MZX Communicator" row (present on the 16-row sheets)
# right above the "RS800" row.
$croatia.Rows.Item(12).Insert()

# Copy the formatting of the row below onto the freshly inserted row, then
# fill in its text.
$croatia.Range("A13").Copy()
$croatia.Range("A12").PasteSpecial(-4122)
$croatia.Range("A12").Value = "MZX Communicator"

# Market-specific values. B4 is set before B2 so the new shared strings are
# appended to the table in the same order as the source workbook.
$croatia.Range("B4").Value = "NGC-3193/T2486/T2485/T2487"
$croatia.Range("B2").Value = "Croatia Market"

# Make the new sheet the active tab/selection, matching the saved workbook
# view state.
$croatia.Activate()
$croatia.Range("B6").Select()
